$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in Q2 and R2 down to plain integers
$ws.Range("Q2").Value = 574678
$ws.Range("R2").Value = 6299720

# Clear the Starttid (Z2) and Sluttid (AB2) cells entirely; Slutdatum (AA2) is unchanged
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
